$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("seats")

# Header: place_id -> loadable_place_id
$ws.Range("C1").Value = "loadable_place_id"

# Row 2 (A2=0, B2 stays "admin")
$ws.Range("C2").Value = "Test event 1; Cinema hall, floor 1, upper middle section, row 6, seat 3"
$ws.Range("D2").Value = 45508.49006479167

# Row 3 (A3=1, B3 "bogdan.yakupov@nu.edu.kz" -> "admin")
$ws.Range("B3").Value = "admin"
$ws.Range("C3").Value = "Test event 1; Cinema hall, floor 1, upper middle section, row 4, seat 2"
$ws.Range("D3").Value = 45508.49010777778

# Row 4 (A4=2, B4 stays "admin")
$ws.Range("C4").Value = "Test event 1; Cinema hall, floor 1, upper middle section, row 8, seat 3"
$ws.Range("D4").Value = 45508.49005549768
